$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.834112048149109
$ws.Range("B1").Value = 3.097206592559814
$ws.Range("C1").Value = 3.431394815444946
$ws.Range("D1").Value = 3.773215055465698
$ws.Range("E1").Value = 2.588592529296875
